$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct "Shapely Ross" -> "Shapley Ross" spelling in the three affected cells.
$ws.Range("C5").Value = "LSR is born to Shapley Ross and Catherine Fulkerson in Iowa (September 27, 1838)"
$ws.Range("C6").Value = "Shapley Ross runs to Texas after physical altercation with a lawyer over a runaway slave (sometime in 1838). Catherine Fulkerson follows with rest of family shortly after. "
$ws.Range("C26").Value = "Shapley Ross, father to LSR, dies (September 17, 1889)"

# Update the active cell selection to match the saved view state.
$ws.Activate()
$ws.Range("C8").Select()
